$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("AG2").Value = 0

$ws.Range("E3").Value = 0
$ws.Range("AG3").Value = 0

$ws.Range("Z4").Value = 0
$ws.Range("AA4").Value = 0

$ws.Range("Z17").Value = 0
$ws.Range("AA17").Value = 0

$ws.Range("AG18").Value = 0

$ws.Range("AG19").Value = 0

$ws.Range("S28").Select()
